$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New vocabulary rows (animal sounds / verbs), appended after the existing 25 data rows.
# Row 27 is intentionally left blank (matches source spreadsheet layout).

$ws.Range("A28").Value = 'roucouler'
$ws.Range("B28").Value = 'vi'
$ws.Range("C28").Value = 'rukule.'
$ws.Range("D28").Value = 'cukrovat, vrkat'

$ws.Range("A29").Value = 'pépier'
$ws.Range("B29").Value = 'vi'
$ws.Range("C29").Value = 'pe.pje.'
$ws.Range("D29").Value = 'pípat, štěbetat'

$ws.Range("A30").Value = 'caqueter'
$ws.Range("B30").Value = 'vi'
$ws.Range("C30").Value = 'kakte.'
$ws.Range("D30").Value = 'kdákat'

$ws.Range("A31").Value = 'jacasser'
$ws.Range("B31").Value = 'vi'
$ws.Range("C31").Value = 'žakase.'
$ws.Range("D31").Value = 'štěbetat'

$ws.Range("A32").Value = 'blatérer'
$ws.Range("B32").Value = 'vi'
$ws.Range("C32").Value = 'blate.re.'
$ws.Range("D32").Value = 'bečet (o beranu), mečet (o velbloudu)'

$ws.Range("A33").Value = '''''huer'
$ws.Range("B33").Value = 'vi'
$ws.Range("C33").Value = 'üe.'
$ws.Range("D33").Value = 'houkat (o sově)'

$ws.Range("A34").Value = 'râler'
$ws.Range("B34").Value = 'vi'
$ws.Range("C34").Value = 'rale.'
$ws.Range("D34").Value = 'chroptět, chrčet, řvát (tygr)'

$ws.Range("A35").Value = 'craquer'
$ws.Range("B35").Value = 'vi'
$ws.Range("C35").Value = 'krake.'
$ws.Range("D35").Value = 'praskat, skřípat, vrzat'

$ws.Range("A36").Value = 'baréter'
$ws.Range("B36").Value = 'vi'
$ws.Range("C36").Value = 'bare.te.'
$ws.Range("D36").Value = 'troubit (o slonu)'

$ws.Range("A37").Value = 'hennir'
$ws.Range("B37").Value = 'vi'
$ws.Range("C37").Value = 'eni:r'
$ws.Range("D37").Value = 'řehtat, ržát (o koni)'

$ws.Range("A38").Value = '''''hululer'
$ws.Range("B38").Value = 'vi'
$ws.Range("C38").Value = 'ülüle.'
$ws.Range("D38").Value = 'houkat (o sově)'

$ws.Range("A39").Value = 'coasser'
$ws.Range("B39").Value = 'vi'
$ws.Range("C39").Value = 'koase.'
$ws.Range("D39").Value = 'kvákat, kuňkat'

$ws.Range("A40").Value = 'piauler'
$ws.Range("B40").Value = 'vi'
$ws.Range("C40").Value = 'pjo.le.'
$ws.Range("D40").Value = 'pípat'

$ws.Range("A41").Value = 'croasser'
$ws.Range("B41").Value = 'vi'
$ws.Range("C41").Value = 'kroase.'
$ws.Range("D41").Value = 'krákorat, krákat (o vráně)'

$ws.Range("A42").Value = 'cajoler'
$ws.Range("B42").Value = 'vi'
$ws.Range("C42").Value = 'kažole.'
$ws.Range("D42").Value = 'křičet (o strace, sojce)'

$ws.Range("A43").Value = 'bramer'
$ws.Range("B43").Value = 'vi'
$ws.Range("C43").Value = 'brame.'
$ws.Range("D43").Value = 'troubit (o jelenu)'

$ws.Range("A44").Value = 'braire'
$ws.Range("B44").Value = 'vi'
$ws.Range("C44").Value = 'bre:r'
$ws.Range("D44").Value = 'hýkat (o oslu)'

$ws.Range("A45").Value = 'cancaner'
$ws.Range("B45").Value = 'vi'
$ws.Range("C45").Value = 'ka~kane.'
$ws.Range("D45").Value = 'káchat (o kachně)'

$ws.Range("A46").Value = 'mugir'
$ws.Range("B46").Value = 'vi'
$ws.Range("C46").Value = 'müži:r'
$ws.Range("D46").Value = 'bučet'

$ws.Range("A47").Value = 'glouglouter'
$ws.Range("B47").Value = 'vi'
$ws.Range("C47").Value = 'gluglute.'
$ws.Range("D47").Value = 'hudrovat (krocan)'

$ws.Range("A48").Value = 'brailler'
$ws.Range("B48").Value = 'vi'
$ws.Range("C48").Value = 'braje.'
$ws.Range("D48").Value = 'řvát, vřeštět (páv)'

$ws.Range("A49").Value = 'margoter'
$ws.Range("B49").Value = 'vi'
$ws.Range("C49").Value = 'margote.'
$ws.Range("D49").Value = 'volat (o křepelce)'

# Re-fill the JSON-building formula in column F down through the new rows so it
# becomes one shared formula group (matches the fill-down the author performed).
$ws.Range("F2:F49").Formula = '= "{ ""foreign"": """ & A2 & """, ""grammar"": """ & B2 & """, ""pronunciation"": """ & C2 & """, ""meaning"": """ & D2 & """ },"'

# Row 27 has no source data (A27:D27 are blank) and is not present in the target
# sheet, so drop the formula cell that the fill created there.
$ws.Range("F27").ClearContents()

